# Atualizei dados bibi e add
# Updates the quarterly recurrence metrics for the last row (2025Q3, row 21)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C21").Value = 297
$ws.Range("D21").Value = 251
$ws.Range("E21").Value = 46
$ws.Range("F21").Value = 71.91977077363897
